$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.866.39'
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.228.61'
$ws.Range("E3").Value = '  -5.00%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '293.17'
$ws.Range("E5").Value = '  -5.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '84.76'
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.512'
$ws.Range("E7").Value = '  -2.64%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -2.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0796'
$ws.Range("E10").Value = '  -1.17%  '
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.96'
$ws.Range("E12").Value = '  -8.70%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.576.00'
$ws.Range("E14").Value = '  -4.88%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.34'
$ws.Range("E15").Value = '  -1.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.15'
$ws.Range("E16").Value = '  -4.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.233.64'
$ws.Range("E17").Value = '  -5.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.722'
$ws.Range("E18").Value = '  -4.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '39.785.77'
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("E20").Value = '  -1.16%  '
$ws.Range("E21").Value = '  -5.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.40'
$ws.Range("E22").Value = '  -4.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.52'
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '232.26'
$ws.Range("E24").Value = '  -1.15%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E26").Value = '  -5.19%  '
$ws.Range("E27").Value = '  +1.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.90'
$ws.Range("E28").Value = '  -3.13%  '
$ws.Range("E29").Value = '  +2.84%  '
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '154.51'
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.75'
$ws.Range("E32").Value = '  -5.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("E34").Value = '  -5.11%  '
$ws.Range("E35").Value = '  -1.73%  '
$ws.Range("E36").Value = '  -5.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.38'
$ws.Range("E37").Value = '  +4.77%  '
$ws.Range("E38").Value = '  -1.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0979'
$ws.Range("E39").Value = '  -1.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.66'
$ws.Range("E40").Value = '  -4.32%  '
$ws.Range("E41").Value = '  -3.70%  '
$ws.Range("E42").Value = '  -3.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.950.44'
$ws.Range("E43").Value = '  -0.88%  '
$ws.Range("E44").Value = '  -2.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0269'
$ws.Range("E45").Value = '  +1.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.36'
$ws.Range("E46").Value = '  -1.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.32'
$ws.Range("E47").Value = '  -6.40%  '
$ws.Range("E48").Value = '  -2.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.444.64'
$ws.Range("E49").Value = '  -4.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.81'
$ws.Range("E50").Value = '  +0.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '88.75'
$ws.Range("E51").Value = '  -4.71%  '
